# Scheduled market-data refresh: update leve crafting profit figures
# (currentAveragePrice / NQ / HQ price & profit columns, H:N) across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets with freshly
# pulled prices. A few rows' HQ price feed dried up (H/I/K -> 0, LeveProfitNQ
# cell removed) or newly appeared (LeveProfitHQ cell added) between runs.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 776110.1
$ws.Range("I62").Value = 940322
$ws.Range("K62").Value = 940322
$ws.Range("M62").Value = -939698
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
# Row 65
$ws.Range("H65").Value = 776110.1
$ws.Range("I65").Value = 940322
$ws.Range("K65").Value = 4701610
$ws.Range("M65").Value = -4698490
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
# Row 74
$ws.Range("H74").Value = 4500
$ws.Range("I74").Value = 4000
$ws.Range("K74").Value = 4000
$ws.Range("M74").Value = -3064
# Row 77
$ws.Range("H77").Value = 4500
$ws.Range("I77").Value = 4000
$ws.Range("K77").Value = 20000
$ws.Range("M77").Value = -15320
# Row 132
$ws.Range("H132").Value = 2634.9
$ws.Range("I132").Value = 2283.194
$ws.Range("K132").Value = 6849.582
$ws.Range("M132").Value = -4319.582
# Row 138
$ws.Range("H138").Value = 1935.9822
$ws.Range("I138").Value = 1498.7073
$ws.Range("K138").Value = 4496.1219
$ws.Range("M138").Value = 643.8780999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 21399.4
$ws.Range("I32").Value = 20499.25
$ws.Range("K32").Value = 20499.25
$ws.Range("M32").Value = -20212.25
# Row 74
$ws.Range("H74").Value = 2821.2126
$ws.Range("I74").Value = 2881.634
$ws.Range("J74").Value = 2408.3333
$ws.Range("K74").Value = 2881.634
$ws.Range("L74").Value = 2408.3333
$ws.Range("M74").Value = -2007.634
$ws.Range("N74").Value = -4156.3333
# Row 77
$ws.Range("H77").Value = 2821.2126
$ws.Range("I77").Value = 2881.634
$ws.Range("J77").Value = 2408.3333
$ws.Range("K77").Value = 14408.17
$ws.Range("L77").Value = 12041.6665
$ws.Range("M77").Value = -10040.17
$ws.Range("N77").Value = -20777.6665
# Row 102
$ws.Range("H102").Value = 47937.59
$ws.Range("I102").Value = 4481.357
$ws.Range("J102").Value = 250733.33
$ws.Range("K102").Value = 4481.357
$ws.Range("L102").Value = 250733.33
$ws.Range("M102").Value = -2859.357
$ws.Range("N102").Value = -253977.33
# Row 122
$ws.Range("H122").Value = 3316.5
$ws.Range("I122").Value = 2999.875
$ws.Range("J122").Value = 3949.75
$ws.Range("K122").Value = 8999.625
$ws.Range("L122").Value = 11849.25
$ws.Range("M122").Value = -6549.625
$ws.Range("N122").Value = -16749.25
# Row 132
$ws.Range("H132").Value = 1249
$ws.Range("I132").Value = 1262.9565
$ws.Range("K132").Value = 3788.8695
$ws.Range("M132").Value = -1258.8695

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 18064.545
$ws.Range("I20").Value = 17076.295
$ws.Range("K20").Value = 17076.295
$ws.Range("M20").Value = -16829.295
# Row 80
$ws.Range("H80").Value = 279.875
$ws.Range("I80").Value = 260
$ws.Range("J80").Value = 339.5
$ws.Range("K80").Value = 260
$ws.Range("L80").Value = 339.5
$ws.Range("M80").Value = 738
$ws.Range("N80").Value = -2335.5
# Row 83
$ws.Range("H83").Value = 279.875
$ws.Range("I83").Value = 260
$ws.Range("J83").Value = 339.5
$ws.Range("K83").Value = 1300
$ws.Range("L83").Value = 1697.5
$ws.Range("M83").Value = 3692
$ws.Range("N83").Value = -11681.5
# Row 99
$ws.Range("H99").Value = 2222.5
$ws.Range("I99").Value = 2284.0527
$ws.Range("K99").Value = 2284.0527
$ws.Range("M99").Value = -786.0527000000002

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5122.151
$ws.Range("I31").Value = 9270.733
$ws.Range("J31").Value = 3484.5527
$ws.Range("K31").Value = 9270.733
$ws.Range("L31").Value = 3484.5527
$ws.Range("M31").Value = -8975.733
$ws.Range("N31").Value = -4074.5527
# Row 34
$ws.Range("H34").Value = 5122.151
$ws.Range("I34").Value = 9270.733
$ws.Range("J34").Value = 3484.5527
$ws.Range("K34").Value = 9270.733
$ws.Range("L34").Value = 3484.5527
$ws.Range("M34").Value = -9068.733
$ws.Range("N34").Value = -3888.5527
# Row 122
$ws.Range("H122").Value = 3011.353
$ws.Range("I122").Value = 3484.3333
$ws.Range("K122").Value = 10452.9999
$ws.Range("M122").Value = -8002.999899999999
# Row 132
$ws.Range("H132").Value = 1884.3556
$ws.Range("I132").Value = 1829.2195
$ws.Range("K132").Value = 5487.6585
$ws.Range("M132").Value = -2957.6585
# Row 134
$ws.Range("H134").Value = 1242
$ws.Range("I134").Value = 1072.0303
$ws.Range("K134").Value = 3216.0909
$ws.Range("M134").Value = -681.0908999999997

$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 500.66666
$ws.Range("I36").Value = 167.33333
$ws.Range("J36").Value = 834
$ws.Range("K36").Value = 501.99999
$ws.Range("L36").Value = 2502
$ws.Range("M36").Value = -332.99999
$ws.Range("N36").Value = -2840
# Row 38
$ws.Range("H38").Value = 100.4375
$ws.Range("I38").Value = 107.64286
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 322.92858
$ws.Range("L38").Value = 150
$ws.Range("M38").Value = 24.07141999999999
$ws.Range("N38").Value = -844

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1206.3793
$ws.Range("I97").Value = 1013.381
$ws.Range("J97").Value = 1713
$ws.Range("K97").Value = 1013.381
$ws.Range("L97").Value = 1713
$ws.Range("M97").Value = -517.381
$ws.Range("N97").Value = -2705
# Row 113
$ws.Range("H113").Value = 1970.7333
$ws.Range("J113").Value = 3020
$ws.Range("L113").Value = 3020
$ws.Range("N113").Value = -7360
# Row 122
$ws.Range("H122").Value = 2610.3845
$ws.Range("I122").Value = 2357.7273
$ws.Range("K122").Value = 7073.1819
$ws.Range("M122").Value = -4623.1819
# Row 126
$ws.Range("H126").Value = 6844.852
$ws.Range("I126").Value = 9507.929
$ws.Range("K126").Value = 28523.787
$ws.Range("M126").Value = -26053.787
# Row 132
$ws.Range("H132").Value = 3283.2
$ws.Range("I132").Value = 2750.0356
$ws.Range("K132").Value = 8250.106800000001
$ws.Range("M132").Value = -5720.106800000001
# Row 133
$ws.Range("H133").Value = 60919.6
$ws.Range("J133").Value = 61275
$ws.Range("L133").Value = 61275
$ws.Range("N133").Value = -71395

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2936888
$ws.Range("J22").Value = 5383320
$ws.Range("L22").Value = 5383320
$ws.Range("N22").Value = -5383910
# Row 27
$ws.Range("H27").Value = 2936888
$ws.Range("J27").Value = 5383320
$ws.Range("L27").Value = 5383320
$ws.Range("N27").Value = -5383534
# Row 59
$ws.Range("H59").Value = 40374.5
$ws.Range("J59").Value = 39999
$ws.Range("L59").Value = 39999
$ws.Range("N59").Value = -41307
# Row 61
$ws.Range("H61").Value = 1406.7391
$ws.Range("I61").Value = 1469.0476
$ws.Range("J61").Value = 752.5
$ws.Range("K61").Value = 1469.0476
$ws.Range("L61").Value = 752.5
$ws.Range("M61").Value = -1267.0476
$ws.Range("N61").Value = -1156.5
# Row 82
$ws.Range("H82").Value = 3088.7273
$ws.Range("I82").Value = 2545.6924
$ws.Range("J82").Value = 3873.111
$ws.Range("K82").Value = 2545.6924
$ws.Range("L82").Value = 3873.111
$ws.Range("M82").Value = -2184.6924
$ws.Range("N82").Value = -4595.111
# Row 85
$ws.Range("H85").Value = 3088.7273
$ws.Range("I85").Value = 2545.6924
$ws.Range("J85").Value = 3873.111
$ws.Range("K85").Value = 2545.6924
$ws.Range("L85").Value = 3873.111
$ws.Range("M85").Value = -1297.6924
$ws.Range("N85").Value = -6369.111
# Row 113
$ws.Range("H113").Value = 1406.7391
$ws.Range("I113").Value = 1469.0476
$ws.Range("J113").Value = 752.5
$ws.Range("K113").Value = 1469.0476
$ws.Range("L113").Value = 752.5
$ws.Range("M113").Value = 700.9523999999999
$ws.Range("N113").Value = -5092.5

$ws = $wb.Worksheets.Item("WVR")
# Row 27
$ws.Range("H27").Value = 49950
$ws.Range("J27").Value = 49950
$ws.Range("L27").Value = 49950
$ws.Range("N27").Value = -50088
# Row 109
$ws.Range("H109").Value = 65687.5
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
# Row 132
$ws.Range("H132").Value = 3623.359
$ws.Range("I132").Value = 3108.8333
$ws.Range("J132").Value = 5338.4443
$ws.Range("K132").Value = 9326.499899999999
$ws.Range("L132").Value = 16015.3329
$ws.Range("M132").Value = -6796.499899999999
$ws.Range("N132").Value = -21075.3329
